$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing phone number in A2 to an invalid (too long) number
$ws.Range("A2").Value = 4798802111777

# Add new row 3 with the original (valid) phone number and name
$ws.Range("A3").Value = 47988021117
$ws.Range("B3").Value = "Marco"

# Adjust column A width (also clears the old bestFit auto-size flag)
$ws.Columns.Item(1).ColumnWidth = 12.5

# Update selection to reflect the new active cell
$ws.Range("C15").Select()
